$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.833.04'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.812.31'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4968'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3878'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09622'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +21.45%  '

$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.28'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.429'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.05%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.813.22'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.17%  '

$ws.Range("B15").Value = 'Solana'
$ws.Range("C15").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.44'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.245'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001128'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.12'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06594'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.60%  '

$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.13'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.937'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.903.02'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("E25").Value = '  -0.91%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.71'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.20%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.91'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.032.90'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.401'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.03'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1069'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.051'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.578'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.625'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06809'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.96%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.959'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02313'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2145'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.35'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -7.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.931'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.73%  '

$ws.Range("E41").Value = '  +1.03%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.145'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5912'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.294'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.690'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.87'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.952'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.176'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06793'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.25%  '
